$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions): update 想去人数 (F column) values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 6282
$wsExhibit.Range("F8").Value = 1384

# Sheet "全部类型" (all types): update 想去人数 (F column) values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 6282
$wsAll.Range("F12").Value = 1384
